# Refresh the scrapedAt (H) and lastSeenAt (I) timestamps for every data
# row (2-46) to the values produced by the latest scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "2026-02-03T14:06:34.992Z"
$ws.Range("I2").Value = "2026-02-03T14:06:35.036Z"
$ws.Range("H3").Value = "2026-02-03T14:06:34.993Z"
$ws.Range("I3").Value = "2026-02-03T14:06:35.059Z"
$ws.Range("H4").Value = "2026-02-03T14:06:34.993Z"
$ws.Range("I4").Value = "2026-02-03T14:06:35.063Z"
$ws.Range("H5").Value = "2026-02-03T14:06:34.993Z"
$ws.Range("I5").Value = "2026-02-03T14:06:35.065Z"
$ws.Range("H6").Value = "2026-02-03T14:06:34.993Z"
$ws.Range("I6").Value = "2026-02-03T14:06:35.067Z"
$ws.Range("H7").Value = "2026-02-03T14:06:34.993Z"
$ws.Range("I7").Value = "2026-02-03T14:06:35.072Z"
$ws.Range("H8").Value = "2026-02-03T14:06:34.994Z"
$ws.Range("I8").Value = "2026-02-03T14:06:35.076Z"
$ws.Range("H9").Value = "2026-02-03T14:06:34.994Z"
$ws.Range("I9").Value = "2026-02-03T14:06:35.079Z"
$ws.Range("H10").Value = "2026-02-03T14:06:34.994Z"
$ws.Range("I10").Value = "2026-02-03T14:06:35.083Z"
$ws.Range("H11").Value = "2026-02-03T14:06:34.994Z"
$ws.Range("I11").Value = "2026-02-03T14:06:35.086Z"
$ws.Range("H12").Value = "2026-02-03T14:06:34.994Z"
$ws.Range("I12").Value = "2026-02-03T14:06:35.090Z"
$ws.Range("H13").Value = "2026-02-03T14:06:34.994Z"
$ws.Range("I13").Value = "2026-02-03T14:06:35.093Z"
$ws.Range("H14").Value = "2026-02-03T14:06:34.995Z"
$ws.Range("I14").Value = "2026-02-03T14:06:35.100Z"
$ws.Range("H15").Value = "2026-02-03T14:06:34.995Z"
$ws.Range("I15").Value = "2026-02-03T14:06:35.104Z"
$ws.Range("H16").Value = "2026-02-03T14:06:34.995Z"
$ws.Range("I16").Value = "2026-02-03T14:06:35.108Z"
$ws.Range("H17").Value = "2026-02-03T14:06:34.995Z"
$ws.Range("I17").Value = "2026-02-03T14:06:35.110Z"
$ws.Range("H18").Value = "2026-02-03T14:06:34.995Z"
$ws.Range("I18").Value = "2026-02-03T14:06:35.113Z"
$ws.Range("H19").Value = "2026-02-03T14:06:34.995Z"
$ws.Range("I19").Value = "2026-02-03T14:06:35.115Z"
$ws.Range("H20").Value = "2026-02-03T14:06:34.995Z"
$ws.Range("I20").Value = "2026-02-03T14:06:35.118Z"
$ws.Range("H21").Value = "2026-02-03T14:06:34.996Z"
$ws.Range("I21").Value = "2026-02-03T14:06:35.120Z"
$ws.Range("H22").Value = "2026-02-03T14:06:34.996Z"
$ws.Range("I22").Value = "2026-02-03T14:06:35.122Z"
$ws.Range("H23").Value = "2026-02-03T14:06:34.996Z"
$ws.Range("I23").Value = "2026-02-03T14:06:35.125Z"
$ws.Range("H24").Value = "2026-02-03T14:06:34.996Z"
$ws.Range("I24").Value = "2026-02-03T14:06:35.128Z"
$ws.Range("H25").Value = "2026-02-03T14:06:34.996Z"
$ws.Range("I25").Value = "2026-02-03T14:06:35.130Z"
$ws.Range("H26").Value = "2026-02-03T14:06:34.996Z"
$ws.Range("I26").Value = "2026-02-03T14:06:35.133Z"
$ws.Range("H27").Value = "2026-02-03T14:06:34.997Z"
$ws.Range("I27").Value = "2026-02-03T14:06:35.136Z"
$ws.Range("H28").Value = "2026-02-03T14:06:34.997Z"
$ws.Range("I28").Value = "2026-02-03T14:06:35.141Z"
$ws.Range("H29").Value = "2026-02-03T14:06:34.997Z"
$ws.Range("I29").Value = "2026-02-03T14:06:35.143Z"
$ws.Range("H30").Value = "2026-02-03T14:06:34.997Z"
$ws.Range("I30").Value = "2026-02-03T14:06:35.146Z"
$ws.Range("H31").Value = "2026-02-03T14:06:34.997Z"
$ws.Range("I31").Value = "2026-02-03T14:06:35.148Z"
$ws.Range("H32").Value = "2026-02-03T14:06:34.997Z"
$ws.Range("I32").Value = "2026-02-03T14:06:35.151Z"
$ws.Range("H33").Value = "2026-02-03T14:06:34.997Z"
$ws.Range("I33").Value = "2026-02-03T14:06:35.153Z"
$ws.Range("H34").Value = "2026-02-03T14:06:34.997Z"
$ws.Range("I34").Value = "2026-02-03T14:06:35.160Z"
$ws.Range("H35").Value = "2026-02-03T14:06:34.998Z"
$ws.Range("I35").Value = "2026-02-03T14:06:35.162Z"
$ws.Range("H36").Value = "2026-02-03T14:06:34.998Z"
$ws.Range("I36").Value = "2026-02-03T14:06:35.164Z"
$ws.Range("H37").Value = "2026-02-03T14:06:34.998Z"
$ws.Range("I37").Value = "2026-02-03T14:06:35.167Z"
$ws.Range("H38").Value = "2026-02-03T14:06:34.999Z"
$ws.Range("I38").Value = "2026-02-03T14:06:35.170Z"
$ws.Range("H39").Value = "2026-02-03T14:06:34.999Z"
$ws.Range("I39").Value = "2026-02-03T14:06:35.173Z"
$ws.Range("H40").Value = "2026-02-03T14:06:34.999Z"
$ws.Range("I40").Value = "2026-02-03T14:06:35.176Z"
$ws.Range("H41").Value = "2026-02-03T14:06:34.999Z"
$ws.Range("I41").Value = "2026-02-03T14:06:35.178Z"
$ws.Range("H42").Value = "2026-02-03T14:06:35.000Z"
$ws.Range("I42").Value = "2026-02-03T14:06:35.180Z"
$ws.Range("H43").Value = "2026-02-03T14:06:35.000Z"
$ws.Range("I43").Value = "2026-02-03T14:06:35.183Z"
$ws.Range("H44").Value = "2026-02-03T14:06:35.001Z"
$ws.Range("I44").Value = "2026-02-03T14:06:35.185Z"
$ws.Range("H45").Value = "2026-02-03T14:06:35.001Z"
$ws.Range("I45").Value = "2026-02-03T14:06:35.188Z"
$ws.Range("H46").Value = "2026-02-03T14:06:35.001Z"
$ws.Range("I46").Value = "2026-02-03T14:06:35.190Z"
